# Applies the "automatic update" edit described in the commit diff:
#  - Column C ("Förändrad") on all data rows (2-14) bumps from 46064 -> 46065
#  - Rows 7-12 get their A (Beteckning), B (Datum) and G (Area ha) values
#    re-shuffled among themselves (everything else on those rows is untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bump "Förändrad" (column C) for every data row from 46064 to 46065 ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46065
}

# --- 2. Re-shuffle Beteckning / Datum / Area for rows 7-12 ---
$newRowData = @{
    7  = @("A 19922-2025", 45771.63034722222, 10.1)
    8  = @("A 25015-2023", 45085.6989699074, 1.8)
    9  = @("A 28266-2025", 45818.56381944445, 1.9)
    10 = @("A 62884-2021", 44504, 0.8)
    11 = @("A 14271-2021", 44278, 6.7)
    12 = @("A 25634-2025", 45803.59570601852, 6)
}

foreach ($r in $newRowData.Keys) {
    $vals = $newRowData[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 7).Value = $vals[2]
}
